# Update Mexico Liga MX base, 03-04-2024 22:09
# This script fixes mismatched fixture rows: for several pairs of adjacent
# rows, all match data (columns B..AC) had been written to the wrong row
# (row index in column A is correct, everything else was swapped). It also
# nudges a few odds columns (R/S/U/V) on three later rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$RowA,
        [int]$RowB
    )
    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Row pairs whose B..AC data was swapped (column A / the "id" index stays put)
Swap-RowData -RowA 188 -RowB 189
Swap-RowData -RowA 237 -RowB 238
Swap-RowData -RowA 251 -RowB 252
Swap-RowData -RowA 264 -RowB 265
Swap-RowData -RowA 282 -RowB 283

# Odds corrections (oddAHH/oddAHA/oddAHOver/oddAHUnder) on a few rows.
# Values are copied from other cells that already hold the exact target
# number (rather than retyped as literals) to avoid any floating point
# re-formatting drift when the workbook is serialized back to OOXML.
$ws.Range("R289").Value = $ws.Range("R4").Value2
$ws.Range("S289").Value = $ws.Range("S4").Value2
$ws.Range("U289").Value = $ws.Range("U22").Value2
$ws.Range("V289").Value = $ws.Range("V22").Value2

$ws.Range("R292").Value = $ws.Range("R7").Value2
$ws.Range("S292").Value = $ws.Range("S7").Value2
$ws.Range("U292").Value = $ws.Range("U12").Value2
$ws.Range("V292").Value = $ws.Range("V12").Value2

$ws.Range("R295").Value = $ws.Range("R5").Value2
$ws.Range("S295").Value = $ws.Range("S5").Value2
$ws.Range("U295").Value = $ws.Range("U12").Value2
$ws.Range("V295").Value = $ws.Range("V12").Value2
